# "Generate Report for Handoff"
# The localization pipeline re-ran: the working file's GUID-based name
# changed, new content-hash XLIFF file names were produced for each target
# locale, and the handoff/generation timestamps advanced.

$wb = $excel.ActiveWorkbook

$newGuidName = "5d4c8901-f4e3-4ce9-ae4a-8f32db70d649.md"

$newPathAndName = "e2e\5d4c8901-f4e3-4ce9-ae4a-8f32db70d649.md"

$newOverviewDate = "2016-10-18 04:33:43"

$newZhCnXlf = "5d4c8901-f4e3-4ce9-ae4a-8f32db70d649.72322a147aea4a47833bfcd275e1b687572dc514.zh-cn.xlf"
$newZhCnHandoffDate = "2016-10-18 04:33:28"

$newDeDeXlf = "5d4c8901-f4e3-4ce9-ae4a-8f32db70d649.72322a147aea4a47833bfcd275e1b687572dc514.de-de.xlf"
$newDeDeHandoffDate = "2016-10-18 04:33:43"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("B2").Hyperlinks.Item(1).TextToDisplay = $newPathAndName
$wsOverview.Range("G2").Value = $newOverviewDate

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuidName
$wsZhCn.Range("A2").Hyperlinks.Item(1).TextToDisplay = $newGuidName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuidName
$wsDeDe.Range("A2").Hyperlinks.Item(1).TextToDisplay = $newGuidName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newDeDeHandoffDate
